# https://jira.hl7.org/browse/FHIR-36285
# Standardize on "CRD Client" rather than EMR/EHR and clarify multi-part
# client wording; also refresh the (stale) cached "last printed" date
# placeholder that Office stamps into the slide master / layouts.

$p = $ppt.ActivePresentation

$oldDate = "2019-03-07"
$newDate = "2022-11-16"
$ppPlaceholderDate = 16

function Update-DatePlaceholder {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if (-not $shp.HasTextFrame) { continue }

        $isDatePlaceholder = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }

        if (-not $isDatePlaceholder) { continue }

        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq $oldDate) {
            $tr.Text = $newDate
        }
    }
}

# 1) Slide master date placeholder.
Update-DatePlaceholder $p.SlideMaster.Shapes

# 2) Every slide layout's date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# 3) Wording fix on slide 1: "EMR action" -> "system action" (CRD Client
#    standardization, not EMR/EHR-specific).
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shp = $slide.Shapes.Item($shi)
        if (-not $shp.HasTextFrame) { continue }
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "2a. Provider performs EMR action") {
            $tr.Text = "2a. Provider performs system action"
        }
    }
}
